$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: becomes a copy of the hyperlinked "username" row, with its own new hyperlink ---
$ws.Range("A3").Value = "sasikala.ars@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:sasikala.ars@gmail.com")
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").Value = "jhj"

# --- Row 4: the text changes but the mailto hyperlink on A4 is kept ---
$ws.Range("A4").Value = "ewrfethyg"
$ws.Range("A4").Style = "Normal"

# --- Row 6: brand-new row of plain (non-hyperlinked) data ---
$ws.Range("A6").Value = "sasikala.ars@gmail.com"
$ws.Range("B6").Value = "tyhjyjthgjnghn"

# --- Apply the new bordered / non-hyperlink "normal" look to every data cell
#     that needs it (B2:B3, A4:B6) ---
foreach ($addr in @("B2", "B3", "A4", "B4", "A5", "B5", "A6", "B6")) {
    $c = $ws.Range($addr)
    $c.Borders.LineStyle = 1
    $c.Borders.Color = 0
}

# --- Column A got a bit narrower ---
$ws.Columns.Item(1).ColumnWidth = 23.3

$ws.Range("B6").Select() | Out-Null
